# The deck's theme (ppt/theme/theme1.xml, used by the one-and-only
# slide master) is being swapped from the custom "Integral / Red Violet"
# palette to the stock Office "Office Theme" palette.
#
# PowerPoint's COM object model doesn't let a script overwrite a whole
# theme part wholesale; the supported, persisting lever is the
# ThemeColorScheme collection (12 theme colors: dk1/lt1/dk2/lt2,
# accent1-6, hlink, folHlink), which is exactly what changed between
# the two <a:clrScheme> blocks in the diff (fonts/format scheme are
# untouched). Updating it from any slide updates the shared master
# theme for the whole presentation.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# .RGB takes the standard OLE color long: R + G*256 + B*65536.
$tcs.Item(1).RGB  = 0x000000   # dk1      000000
$tcs.Item(2).RGB  = 0xFFFFFF   # lt1      FFFFFF
$tcs.Item(3).RGB  = 0x6A5444   # dk2      44546A
$tcs.Item(4).RGB  = 0xE6E6E7   # lt2      E7E6E6
$tcs.Item(5).RGB  = 0xD59B5B   # accent1  5B9BD5
$tcs.Item(6).RGB  = 0x317DED   # accent2  ED7D31
$tcs.Item(7).RGB  = 0xA5A5A5   # accent3  A5A5A5
$tcs.Item(8).RGB  = 0x00C0FF   # accent4  FFC000
$tcs.Item(9).RGB  = 0xC47244   # accent5  4472C4
$tcs.Item(10).RGB = 0x47AD70   # accent6  70AD47
$tcs.Item(11).RGB = 0xC16305   # hlink    0563C1
$tcs.Item(12).RGB = 0x724F95   # folHlink 954F72
